# Append one new telemetry row to the bottom of each of the 4 worksheets
# (file-upload refresh of the SAG2 database).

function AddDataRow($ws, $row, $aVal, $bVal, $cVal, $dVal, $eVal, $fVal, $gVal, $hVal, $iVal) {
    # Column A: date/time serial. Re-apply the same custom number format used
    # by the rest of the column so the new cell reuses the existing style.
    $ws.Cells.Item($row, 1).Value = $aVal
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = $bVal
    $ws.Cells.Item($row, 3).Value = $cVal
    $ws.Cells.Item($row, 4).Value = $dVal
    $ws.Cells.Item($row, 5).Value = $eVal
    $ws.Cells.Item($row, 6).Value = $fVal
    $ws.Cells.Item($row, 7).Value = $gVal
    $ws.Cells.Item($row, 8).Value = $hVal
    $ws.Cells.Item($row, 9).Value = $iVal
}

$wb = $excel.ActiveWorkbook

# Large ID_DEC values overflow double precision round-tripping, so build the
# numeric one explicitly from its string form (this parser doesn't accept
# bare scientific-notation literals like 5.68631262647114e+23).
$idDecNumeric = [double]"5.68631262647114e+23"

# --- Sheet 1: ROW50-FE-LIFTER -> new row 18 ---
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
AddDataRow $ws1 18 45733.11858277778 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x82" "0xe" 400 $idDecNumeric 386 14

# --- Sheet 2: ROW50-MID-LIFTER -> new row 20 ---
# This sheet's ID_DEC column (G) is stored as text throughout (the value is
# too large to round-trip exactly as a double), and column B keeps a
# trailing space like every other row on this sheet. Leading "'" forces the
# big numeric string to be kept as text instead of being parsed as a number.
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
AddDataRow $ws2 20 45733.09563657407 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x86" "0x19" 400 "'568631262647113771663628" 390 25

# --- Sheet 3: ROW11-FE-LIFTER -> new row 18 ---
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
AddDataRow $ws3 18 45733.13871136574 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x82" "0x14" 400 $idDecNumeric 386 20

# --- Sheet 4: ROW11-MID-LIFTER -> new row 18 ---
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
AddDataRow $ws4 18 45733.2857841088 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x86" "0x19" 400 $idDecNumeric 390 25
